$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text storage
# (COM.Value would otherwise auto-infer these digit/dot strings as numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.034.96"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").Value = "1.827.54"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.51%  "

$ws.Range("D5").Value = "311.84"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").Value = "0.4347"
$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("D8").Value = "0.3673"
$ws.Range("E8").Value = "  -1.17%  "

$ws.Range("D9").Value = "0.07311"
$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").Value = "0.8447"
$ws.Range("E10").Value = "  -2.90%  "

$ws.Range("D11").Value = "20.73"
$ws.Range("E11").Value = "  -2.49%  "

$ws.Range("D12").Value = "1.828.81"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "6.666"
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.301"
$ws.Range("E14").Value = "  -1.55%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.07059"
$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").Value = "89.52"
$ws.Range("E16").Value = "  +1.65%  "

$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("D18").Value = "0.000008794"
$ws.Range("E18").Value = "  -1.60%  "

$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").Value = "14.91"
$ws.Range("E20").Value = "  -2.65%  "

$ws.Range("D21").Value = "27.065.39"
$ws.Range("E21").Value = "  -1.68%  "

$ws.Range("D22").Value = "5.145"
$ws.Range("E22").Value = "  -0.94%  "

$ws.Range("D23").Value = "10.88"
$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("D24").Value = "2.050.69"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("D25").Value = "1.984"
$ws.Range("E25").Value = "  -1.20%  "

$ws.Range("D26").Value = "151.41"
$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").Value = "2.220"
$ws.Range("E27").Value = "  +2.56%  "

$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").Value = "5.242"
$ws.Range("E29").Value = "  -1.01%  "

$ws.Range("D30").Value = "117.02"
$ws.Range("E30").Value = "  -0.82%  "

$ws.Range("D31").Value = "0.08738"
$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("D32").Value = "1.181"
$ws.Range("E32").Value = "  -2.69%  "

$ws.Range("D33").Value = "0.7412"
$ws.Range("E33").Value = "  -3.81%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.909"
$ws.Range("E34").Value = "  -0.09%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "4.442"
$ws.Range("E35").Value = "  -1.41%  "

$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").Value = "1.095"
$ws.Range("E37").Value = "  -2.70%  "

$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("D39").Value = "0.05236"
$ws.Range("E39").Value = "  -1.20%  "

$ws.Range("D40").Value = "7.227"
$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("D42").Value = "0.1704"
$ws.Range("E42").Value = "  +1.55%  "

$ws.Range("D43").Value = "0.5133"
$ws.Range("E43").Value = "  +0.50%  "

$ws.Range("D44").Value = "8.599"
$ws.Range("E44").Value = "  -1.60%  "

$ws.Range("D45").Value = "10.58"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "1.947"
$ws.Range("E47").Value = "  +6.23%  "

$ws.Range("D48").Value = "105.97"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").Value = "0.9998"
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("D50").Value = "1.663"
$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("D51").Value = "0.06338"
$ws.Range("E51").Value = "  -1.50%  "
